$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("is_active") for data rows 2-11 currently stores a =TRUE()
# formula that evaluates to the number 1. The fix replaces each of
# those cells with the literal text "TRUE" instead of the boolean
# formula result.
$rng = $ws.Range("D2:D11")
$rng.NumberFormat = "@"
$rng.Formula = '=TEXT(TRUE(),"")'

$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Reflect the final active selection as captured in the saved workbook.
$ws.Range("F9").Select()
